$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content changes -------------------------------------------------
# A7 used to read "In addition, your task will be to answer job-related
# questions." It now reads a short lead-in line, and the detailed
# sub-bullets that used to live in A8:A13 are removed (cells blanked, but
# kept as part of the same bordered box).
$ws.Range("A7").Value = "These jobs are as follows:"
$ws.Range("A8:A13").ClearContents()

# --- Formatting changes ----------------------------------------------------
# xlEdgeLeft = 7, xlEdgeTop = 8, xlEdgeBottom = 9, xlEdgeRight = 10
$grayLeft   = 0xAAAAAA   # new lighter gray used for the left edge
$grayRight  = 0x3F3F3F   # existing darker gray used for the right edge
$grayTopBot = 0xA5A5A5   # existing mid gray used for top/bottom edges

# A7: top row of the box - keep a top border, drop the bottom border,
# switch the fill/font to the plain (non-bold, light) style used by the
# rest of the box instead of the bold header look it had before.
$rngA7 = $ws.Range("A7")
$rngA7.Font.Bold = $false
$rngA7.Interior.Color = 0xFFFF00
$rngA7.Borders.Item(7).LineStyle = 1
$rngA7.Borders.Item(7).Color = $grayLeft
$rngA7.Borders.Item(10).LineStyle = 1
$rngA7.Borders.Item(10).Color = $grayRight
$rngA7.Borders.Item(8).LineStyle = 1
$rngA7.Borders.Item(8).Color = $grayTopBot
$rngA7.Borders.Item(9).LineStyle = -4142

# A8:A12: middle of the box - no top or bottom border, just the left/right
# edges continuing down, same plain fill/font.
$rngMid = $ws.Range("A8:A12")
$rngMid.Font.Bold = $false
$rngMid.Interior.Color = 0xFFFF00
$rngMid.Borders.Item(7).LineStyle = 1
$rngMid.Borders.Item(7).Color = $grayLeft
$rngMid.Borders.Item(10).LineStyle = 1
$rngMid.Borders.Item(10).Color = $grayRight
$rngMid.Borders.Item(8).LineStyle = -4142
$rngMid.Borders.Item(9).LineStyle = -4142

# A13: bottom of the box - no top border, but a bottom border closes the
# box off, same plain fill/font.
$rngA13 = $ws.Range("A13")
$rngA13.Font.Bold = $false
$rngA13.Interior.Color = 0xFFFF00
$rngA13.Borders.Item(7).LineStyle = 1
$rngA13.Borders.Item(7).Color = $grayLeft
$rngA13.Borders.Item(10).LineStyle = 1
$rngA13.Borders.Item(10).Color = $grayRight
$rngA13.Borders.Item(8).LineStyle = -4142
$rngA13.Borders.Item(9).LineStyle = 1
$rngA13.Borders.Item(9).Color = $grayTopBot
